$wb = $excel.ActiveWorkbook

# --- Rename "Requested quantity" headers on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy header/date cell formatting from the "Weekly Quantity" sheet so the
# new sheet reuses the same styles (bold/bordered header, date number format).
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A42").PasteSpecial(-4122)

# --- Header row ---
$headers = @("ds", "PO_Forecast", "yhat_lower", "yhat_upper")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $wsForecast.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# --- Data rows (41 rows: rows 2-42) ---
$poForecastData = @(
    @(45221.99999999999, 2, -13.14143821145515, 20.10213338325745),
    @(45228.99999999999, 2, -15.10350180492543, 17.74673082536317),
    @(45249.99999999999, 3, -13.10033481898137, 20.10611056512057),
    @(45256.99999999999, 4, -13.22579717545753, 20.49952095567826),
    @(45270.99999999999, 4, -11.88053445404841, 20.68607298504871),
    @(45277.99999999999, 5, -10.79372930881389, 20.9893761043067),
    @(45298.99999999999, 5, -11.26060462549372, 21.75327781153001),
    @(45305.99999999999, 6, -10.06999228405025, 22.09886352023969),
    @(45312.99999999999, 6, -10.39354691968957, 22.50401535404275),
    @(45319.99999999999, 6, -10.14963795330642, 22.28229388882935),
    @(45326.99999999999, 7, -10.44614579470676, 21.48324103476445),
    @(45333.99999999999, 7, -9.484438968987639, 22.42872932805552),
    @(45340.99999999999, 7, -7.66726274253612, 24.4948598629698),
    @(45361.99999999999, 8, -7.351117342376463, 23.93193668398477),
    @(45403.99999999999, 10, -6.560417758942178, 26.89065069878362),
    @(45410.99999999999, 10, -6.496656279609047, 25.77619180104325),
    @(45417.99999999999, 10, -5.726110146044483, 27.57408422123605),
    @(45424.99999999999, 11, -5.563291064174527, 26.51236937954425),
    @(45431.99999999999, 11, -5.068055660370558, 26.87057234639788),
    @(45438.99999999999, 11, -4.672961011733624, 27.17454893236975),
    @(45445.99999999999, 12, -4.795292397850783, 28.13622846586318),
    @(45452.99999999999, 12, -4.740907597617255, 28.91856434932692),
    @(45459.99999999999, 12, -4.587606692672527, 28.16512022864535),
    @(45487.99999999999, 13, -1.965651168576851, 29.45839954937888),
    @(45494.99999999999, 14, -1.602795269469304, 29.69882797301499),
    @(45501.99999999999, 14, -2.106531581069087, 30.95802871039791),
    @(45515.99999999999, 15, -2.413236364676026, 31.15383685396379),
    @(45522.99999999999, 15, -0.5157324871458384, 29.87283266467075),
    @(45529.99999999999, 15, -1.317323771895849, 31.39134547837433),
    @(45536.99999999999, 15, -0.1950227797221908, 31.65189916611168),
    @(45543.99999999999, 16, -0.1356387595893556, 31.49962645735555),
    @(45564.99999999999, 17, 0.7109576285015569, 31.8428178910124),
    @(45571.99999999999, 17, -0.1057900025973702, 33.69660812204293),
    @(45578.99999999999, 17, 0.9910555096696334, 33.65713885082981),
    @(45585.99999999999, 18, 1.487371238978906, 33.27371928688771),
    @(45592.99999999999, 18, 1.418070789668806, 33.67079976290731),
    @(45599.99999999999, 18, 4.228182285563139, 34.55234316011227),
    @(45606.99999999999, 18, 2.026510102872942, 33.45649635260392),
    @(45613.99999999999, 19, 2.155028892598436, 34.81436667391626),
    @(45620.99999999999, 19, 3.964930802039572, 34.91284242603662),
    @(45627.99999999999, 19, 3.69689517641462, 34.71062241372786)
)

for ($r = 0; $r -lt $poForecastData.Count; $r++) {
    $row = $poForecastData[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $wsForecast.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
